$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Japan indices for 2019 (row 25) and 2020 (row 26)
$ws.Range("D25").Value = 101.4
$ws.Range("E25").Value = 100.6
$ws.Range("D26").Value = 100
$ws.Range("E26").Value = 100
